# Auto update Excel log
# Appends new sensor/log rows to four worksheets (ALERTS, Proximity, mmWave, Camera).
# Columns in every sheet: Date | Timestamp | Hour | Location | Value | Status

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALERTS (sheet1) — new rows 13:14
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALERTS")
$data = @(
    ,@("2026-02-01", "14:42:43", "14:00", "Living Room", "CRITICAL", "FALL_DETECTED")
    ,@("2026-02-01", "14:42:45", "14:00", "Living Room", "CRITICAL", "FALL_DETECTED")
)
$r = 13
foreach ($row in $data) {
    # Column A looks like a date ("2026-02-01"); force Text so Excel doesn't
    # silently convert it into a date serial number.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# ---------------------------------------------------------------------------
# Proximity (sheet5) — new rows 34:38
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Proximity")
$data = @(
    ,@("2026-02-01", "14:42:37", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
    ,@("2026-02-01", "14:42:45", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
    ,@("2026-02-01", "14:42:46", "14:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door")
    ,@("2026-02-01", "14:42:56", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
    ,@("2026-02-01", "14:42:58", "14:00", "Living Room Main Door", "EXIT",  "User EXITED Living Room Main Door")
)
$r = 34
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# ---------------------------------------------------------------------------
# mmWave (sheet6) — new rows 24:29
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("mmWave")
$data = @(
    ,@("2026-02-01", "14:42:16", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")
    ,@("2026-02-01", "14:42:27", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")
    ,@("2026-02-01", "14:42:37", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")
    ,@("2026-02-01", "14:42:47", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")
    ,@("2026-02-01", "14:42:50", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")
    ,@("2026-02-01", "14:42:58", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")
)
$r = 24
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# ---------------------------------------------------------------------------
# Camera (sheet7) — new rows 22:24
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Camera")
$data = @(
    ,@("2026-02-01", "14:42:39", "14:00", "Living Room Main Door", "Image Received", "Active")
    ,@("2026-02-01", "14:42:46", "14:00", "Living Room Main Door", "Image Received", "Active")
    ,@("2026-02-01", "14:42:58", "14:00", "Living Room Main Door", "Image Captured", "Active")
)
$r = 22
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}
